$d = $word.ActiveDocument

function Replace-Exact($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replaceText, 2)
}

# Main body: RA number run cleared to two spaces
Replace-Exact " 000112544246 - 3 " "  "

# Main body: bold/bCs "QWR" -> "TERE"
Replace-Exact "QWR" "TERE"

# Header: "QWER" -> "TRE"
Replace-Exact "QWER" "TRE"

# Header: remaining "QWR" occurrence -> "TERE"
Replace-Exact "QWR" "TERE"

# Header: five "Qwer" -> "Tre"
Replace-Exact "Qwer" "Tre"

# Header: three "qwer" -> "tre"
Replace-Exact "qwer" "tre"
